$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("STZ")

# Row 4 - Inventory
$ws.Range("B4").Value = 1291000000.0
$ws.Range("C4").Value = 1377000000.0
$ws.Range("D4").Value = 1328000000.0
$ws.Range("E4").Value = 1333000000.0
$ws.Range("F4").Value = 1374000000.0

# Row 15 - Accounts Payable
$ws.Range("B15").Value = 460000000.0
$ws.Range("C15").Value = 732000000.0
$ws.Range("D15").Value = 651000000.0
$ws.Range("E15").Value = 506000000.0
$ws.Range("F15").Value = 558000000.0

# Row 21 - Long Term Tax Liability (Deferred)
$ws.Range("B21").Value = -1654000000.0
$ws.Range("C21").Value = -2563000000.0
$ws.Range("D21").Value = -2575000000.0
$ws.Range("E21").Value = -2600000000.0
$ws.Range("F21").Value = -1996000000.0
